$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2293
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325

$ws.Range("H58").Value = 518.7143
$ws.Range("J58").Value = 602.6667
$ws.Range("L58").Value = 1808.0001
$ws.Range("N58").Value = -2108.0001

$ws.Range("H64").Value = 3500
$ws.Range("I64").Value = 3500
$ws.Range("K64").Value = 3500
$ws.Range("M64").Value = -3252

$ws.Range("H67").Value = 3500
$ws.Range("I67").Value = 3500
$ws.Range("K67").Value = 3500
$ws.Range("M67").Value = -2642

$ws.Range("H76").Value = 6579.8
$ws.Range("I76").Value = 6724.875
$ws.Range("J76").Value = 5999.5
$ws.Range("K76").Value = 6724.875
$ws.Range("L76").Value = 5999.5
$ws.Range("M76").Value = -6409.875
$ws.Range("N76").Value = -6629.5

$ws.Range("H79").Value = 6579.8
$ws.Range("I79").Value = 6724.875
$ws.Range("J79").Value = 5999.5
$ws.Range("K79").Value = 6724.875
$ws.Range("L79").Value = 5999.5
$ws.Range("M79").Value = -5632.875
$ws.Range("N79").Value = -8183.5

$ws.Range("H94").Value = 4154.7144
$ws.Range("I94").Value = 3847.1667
$ws.Range("K94").Value = 3847.1667
$ws.Range("M94").Value = -3396.1667

$ws.Range("H138").Value = 3966
$ws.Range("I138").Value = 4583
$ws.Range("J138").Value = 3595.8
$ws.Range("K138").Value = 13749
$ws.Range("L138").Value = 10787.4
$ws.Range("M138").Value = -8609
$ws.Range("N138").Value = -21067.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 790.2857
$ws.Range("I6").Value = 790.2857
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 790.2857
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -617.2857
$ws.Range("N6").Value = $null

$ws.Range("H36").Value = 4389.8
$ws.Range("I36").Value = 4389.8
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4389.8
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4043.8
$ws.Range("N36").Value = $null

$ws.Range("H97").Value = 561.875
$ws.Range("I97").Value = 580
$ws.Range("K97").Value = 580
$ws.Range("M97").Value = -84

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null

$ws.Range("H132").Value = 12247.25
$ws.Range("J132").Value = 15996.667
$ws.Range("L132").Value = 47990.001
$ws.Range("N132").Value = -53050.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 162.4
$ws.Range("I22").Value = 162.4
$ws.Range("K22").Value = 162.4
$ws.Range("M22").Value = 10.59999999999999

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null

$ws.Range("H94").Value = 1050.2
$ws.Range("I94").Value = 1050.2
$ws.Range("K94").Value = 1050.2
$ws.Range("M94").Value = -599.2

$ws.Range("H102").Value = 13916.5
$ws.Range("I102").Value = 13916.5
$ws.Range("K102").Value = 13916.5
$ws.Range("M102").Value = -10671.5

$ws.Range("H105").Value = 2136.8
$ws.Range("I105").Value = 2046
$ws.Range("K105").Value = 2046
$ws.Range("M105").Value = -299

$ws.Range("H134").Value = 1960.6
$ws.Range("I134").Value = 1960.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5881.799999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3346.799999999999
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3177.375
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 4069.8333
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 4069.8333
$ws.Range("M3").Value = -387
$ws.Range("N3").Value = -4295.8333

$ws.Range("H7").Value = 248.55556
$ws.Range("I7").Value = 275.25
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 275.25
$ws.Range("L7").Value = 35
$ws.Range("M7").Value = -162.25
$ws.Range("N7").Value = -261

$ws.Range("H57").Value = 49600
$ws.Range("J57").Value = 49600
$ws.Range("L57").Value = 49600
$ws.Range("N57").Value = -50720

$ws.Range("H58").Value = 14500
$ws.Range("I58").Value = 7500
$ws.Range("J58").Value = 18000
$ws.Range("K58").Value = 7500
$ws.Range("L58").Value = 18000
$ws.Range("M58").Value = -7297
$ws.Range("N58").Value = -18406

$ws.Range("H134").Value = 6655.857
$ws.Range("I134").Value = 3300
$ws.Range("J134").Value = 7998.2
$ws.Range("K134").Value = 9900
$ws.Range("L134").Value = 23994.6
$ws.Range("M134").Value = -7365
$ws.Range("N134").Value = -29064.6

$ws.Range("H136").Value = 14500
$ws.Range("I136").Value = 7500
$ws.Range("J136").Value = 18000
$ws.Range("K136").Value = 22500
$ws.Range("L136").Value = 54000
$ws.Range("M136").Value = -19950
$ws.Range("N136").Value = -59100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2247.5
$ws.Range("I39").Value = 2263.3333
$ws.Range("J39").Value = 2200
$ws.Range("K39").Value = 6789.999899999999
$ws.Range("L39").Value = 6600
$ws.Range("M39").Value = -6495.999899999999
$ws.Range("N39").Value = -7188

$ws.Range("H82").Value = 13
$ws.Range("I82").Value = 13
$ws.Range("K82").Value = 39
$ws.Range("M82").Value = 367

$ws.Range("H85").Value = 13
$ws.Range("I85").Value = 13
$ws.Range("K85").Value = 39
$ws.Range("M85").Value = 1365

$ws.Range("H107").Value = 185.5
$ws.Range("I107").Value = 80.666664
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 241.999992
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1678.000008
$ws.Range("N107").Value = -5340

$ws.Range("H113").Value = 1214.9
$ws.Range("I113").Value = 1143.5
$ws.Range("K113").Value = 3430.5
$ws.Range("M113").Value = -1260.5

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null

$ws.Range("H10").Value = 5501
$ws.Range("J10").Value = 999
$ws.Range("L10").Value = 999
$ws.Range("N10").Value = -1337

$ws.Range("H113").Value = 2195.5
$ws.Range("I113").Value = 1832
$ws.Range("K113").Value = 1832
$ws.Range("M113").Value = 338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13000
$ws.Range("J136").Value = 13000
$ws.Range("L136").Value = 39000
$ws.Range("N136").Value = -44100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = 0

$ws.Range("H107").Value = 1119.65
$ws.Range("I107").Value = 1135.2727
$ws.Range("K107").Value = 3405.8181
$ws.Range("M107").Value = -1485.8181

$ws.Range("H132").Value = 13998.5
$ws.Range("I132").Value = 12997
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 38991
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -36461
$ws.Range("N132").Value = -50060
